# Arquivo atualizado em 07/12/2023, 14:20:05.
#
# Column B currently holds Excel date-serial values (formatted with the
# custom "YYYY-MM-DD HH:MM:SS" number format, numFmtId 165) that repeat the
# same 12 years (2010-2021) for each of the three regions (Brasil,
# Nordeste, Sergipe). The edit replaces each of those date cells with a
# plain text "dd/mm/yyyy"-style label, and rescales column C (which held a
# 0-1 fraction) up by a factor of 100 so it reads as a percentage value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 12 years, in the order they appear within each 12-row region block.
$years = @(
    "01/01/2010", "01/01/2011", "01/01/2012", "01/01/2013",
    "01/01/2014", "01/01/2015", "01/01/2016", "01/01/2017",
    "01/01/2018", "01/01/2019", "01/01/2020", "01/01/2021"
)

for ($row = 2; $row -le 37; $row++) {
    $yearIndex = ($row - 2) % 12

    # Column B: replace the date-serial value with a literal text label.
    # The leading apostrophe forces Excel to treat it as text instead of
    # re-parsing "01/01/2010" back into a date serial.
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.Value = "'" + $years[$yearIndex]
    # Drop back to the default "Normal" style so the cell no longer carries
    # the old date number format (or the text quote-prefix style) - matches
    # the target, where these cells carry no explicit style at all.
    $bCell.Style = "Normal"

    # Column C: scale the stored fraction up by 100.
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = $cCell.Value() * 100
}
